$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly batch of price records was reported before the existing
# rows, so insert two blank rows at the top of the data block (row 23,
# right after the existing row 22) and shift the rest of the table down.
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(23).Insert()

# Fill in the two newly inserted rows with the new records.
$ws.Cells.Item(23, 1).Value2 = 5
$ws.Cells.Item(23, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(23, 3).Value2 = "Maule"
$ws.Cells.Item(23, 4).Value2 = 45082
$ws.Cells.Item(23, 5).Value2 = 7
$ws.Cells.Item(23, 6).Value2 = "Fruta"
$ws.Cells.Item(23, 7).Value2 = 100107
$ws.Cells.Item(23, 8).Value2 = "Otros"
$ws.Cells.Item(23, 9).Value2 = 100107001
$ws.Cells.Item(23, 10).Value2 = "Caqui"
$ws.Cells.Item(23, 11).Value2 = "Mankaki"
$ws.Cells.Item(23, 12).Value2 = "Especial"
$ws.Cells.Item(23, 13).Value2 = 200
$ws.Cells.Item(23, 14).Value2 = 12000
$ws.Cells.Item(23, 15).Value2 = 12000
$ws.Cells.Item(23, 16).Value2 = 12000
$ws.Cells.Item(23, 17).Value2 = "`$/caja 12 kilos granel"
$ws.Cells.Item(23, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(23, 19).Value2 = 12000
$ws.Cells.Item(23, 20).Value2 = 1

$ws.Cells.Item(24, 1).Value2 = 5
$ws.Cells.Item(24, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(24, 3).Value2 = "Maule"
$ws.Cells.Item(24, 4).Value2 = 45082
$ws.Cells.Item(24, 5).Value2 = 7
$ws.Cells.Item(24, 6).Value2 = "Fruta"
$ws.Cells.Item(24, 7).Value2 = 100107
$ws.Cells.Item(24, 8).Value2 = "Otros"
$ws.Cells.Item(24, 9).Value2 = 100107001
$ws.Cells.Item(24, 10).Value2 = "Caqui"
$ws.Cells.Item(24, 11).Value2 = "Mankaki"
$ws.Cells.Item(24, 12).Value2 = "Primera"
$ws.Cells.Item(24, 13).Value2 = 250
$ws.Cells.Item(24, 14).Value2 = 10000
$ws.Cells.Item(24, 15).Value2 = 10000
$ws.Cells.Item(24, 16).Value2 = 10000
$ws.Cells.Item(24, 17).Value2 = "`$/caja 12 kilos granel"
$ws.Cells.Item(24, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(24, 19).Value2 = 10000
$ws.Cells.Item(24, 20).Value2 = 1
